$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2001")

# Update pie chart source values (fix a data bug / adjust layout data)
$ws.Range("D2").Value = 0.5
$ws.Range("F2").Value = 0.3
$ws.Range("H2").Value = 0.1

$ws.Range("D3").Value = 0.3
$ws.Range("E3").Value = 0.8
$ws.Range("G3").Value = 0.1

$ws.Range("C5").Value = 0.1
$ws.Range("E5").Value = 0.3

# Move the active selection on this sheet to a single cell. Selecting a
# range requires the sheet to be active, so remember which tab was active
# and restore it afterwards.
$prevActive = $wb.ActiveSheet
$ws.Activate()
$ws.Range("E5").Select()
$prevActive.Activate()
